$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6820405125617981
$ws.Range("B1").Value = 1.414020180702209
$ws.Range("C1").Value = 3.680445671081543
$ws.Range("D1").Value = 3.097216606140137
$ws.Range("E1").Value = 1.73179829120636
